# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker/period detail table (rows 16-36, columns B:G) is reorganized
# from "grouped by worker" order to "grouped by period" order, adding a
# new worker (MENCY MAGDALENA ROSSINI CARMONA) into the rotation so each
# period block now lists all three workers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Worker reference data: DocType, DocNumber, Name, ValorMora (normal), SalarioBasico
$rafael  = @{ Doc = "79152080";   Name = "RAFAEL ANTONIO ABONDANO CAPELLA";     Mora = 52577; MoraUlt = 16619; Salario = 1314431 }
$jenny   = @{ Doc = "1143362836"; Name = "JENNYFER KATERINE CASTELLON CASTRO"; Mora = 55680; MoraUlt = 20416; Salario = 1392000 }
$mency   = @{ Doc = "30763028";   Name = "MENCY MAGDALENA ROSSINI CARMONA";    Mora = 52000; MoraUlt = 19067; Salario = 1300000 }

$workers = @($rafael, $jenny, $mency)
$periods = @("2406", "2407", "2408", "2409", "2410", "2411", "2412")

$row = 16
foreach ($period in $periods) {
    foreach ($w in $workers) {
        $ws.Cells.Item($row, 2).Value = "CC"
        $ws.Cells.Item($row, 3).Value = $w.Doc
        $ws.Cells.Item($row, 4).Value = $w.Name
        $ws.Cells.Item($row, 5).Value = $period
        if ($period -eq "2412") {
            $ws.Cells.Item($row, 6).Value = $w.MoraUlt
        } else {
            $ws.Cells.Item($row, 6).Value = $w.Mora
        }
        $ws.Cells.Item($row, 7).Value = $w.Salario
        $row = $row + 1
    }
}
